$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: time_difference -------------------------------------

# Header cell (G1) - copy the header style from the neighbouring header
# cell (F1) so it keeps the same bold/border/centered look.
$ws.Range("G1").Value = "time_difference"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# Rows 2-4 have no arrival/departure info recorded, so the new column is
# blank for them too - write an empty text value (quote-prefix empty
# string) and then drop back to the default "Normal" style so no stray
# formatting is left behind.
$ws.Range("G2").Formula = "'"
$ws.Range("G2").Style = "Normal"
$ws.Range("G3").Formula = "'"
$ws.Range("G3").Style = "Normal"
$ws.Range("G4").Formula = "'"
$ws.Range("G4").Style = "Normal"

# --- Row 5: corrected arrival/departure date & time ---------------------
# Force text storage (not a date/time serial) by setting the number
# format to Text before assigning, then reset the style back to Normal
# so no extra number-format style is left on the cell.
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "2023-07-25"
$ws.Range("C5").Style = "Normal"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "16:44:47"
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2023-07-25"
$ws.Range("E5").Style = "Normal"

$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "16:44:56"
$ws.Range("F5").Style = "Normal"

# Computed time difference between arrival and departure (departure - arrival).
$ws.Range("G5").Value = "0:00:09"
